# Auto-generated: update cached numeric values in Hyperion_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 820.9524
$ws.Range("I15").Value2 = 820.9524
$ws.Range("K15").Value2 = 2462.8572
$ws.Range("M15").Value2 = -2293.8572
$ws.Range("H28").Value2 = 460.75
$ws.Range("I28").Value2 = 411.81818
$ws.Range("K28").Value2 = 411.81818
$ws.Range("M28").Value2 = 73.18182000000002
$ws.Range("H62").Value2 = 4064.3547
$ws.Range("J62").Value2 = 11153.2
$ws.Range("L62").Value2 = 11153.2
$ws.Range("N62").Value2 = -12401.2
$ws.Range("H65").Value2 = 4064.3547
$ws.Range("J65").Value2 = 11153.2
$ws.Range("L65").Value2 = 55766
$ws.Range("N65").Value2 = -62006
$ws.Range("H82").Value2 = 3997.6667
$ws.Range("I82").Value2 = 3997.6667
$ws.Range("K82").Value2 = 11993.0001
$ws.Range("M82").Value2 = -11587.0001
$ws.Range("H85").Value2 = 3997.6667
$ws.Range("I85").Value2 = 3997.6667
$ws.Range("K85").Value2 = 11993.0001
$ws.Range("M85").Value2 = -10589.0001
$ws.Range("H95").Value2 = 27466.666
$ws.Range("J95").Value2 = 42400
$ws.Range("L95").Value2 = 42400
$ws.Range("N95").Value2 = -47892
$ws.Range("H97").Value2 = 3026.25
$ws.Range("J97").Value2 = 3368.3333
$ws.Range("L97").Value2 = 10104.9999
$ws.Range("N97").Value2 = -11096.9999
$ws.Range("H111").Value2 = 11111655
$ws.Range("I111").Value2 = 11111655
$ws.Range("K111").Value2 = 33334965
$ws.Range("M111").Value2 = -33331898
$ws.Range("H112").Value2 = 6248.091
$ws.Range("J112").Value2 = 7047.579
$ws.Range("L112").Value2 = 21142.737
$ws.Range("N112").Value2 = -23358.737
$ws.Range("H113").Value2 = 5932.2085
$ws.Range("I113").Value2 = 5749.25
$ws.Range("J113").Value2 = 6023.6875
$ws.Range("K113").Value2 = 5749.25
$ws.Range("L113").Value2 = 6023.6875
$ws.Range("M113").Value2 = -2495.25
$ws.Range("N113").Value2 = -12531.6875
$ws.Range("H116").Value2 = 5349.6665
$ws.Range("I116").Value2 = 4661.5
$ws.Range("J116").Value2 = 6037.8335
$ws.Range("K116").Value2 = 4661.5
$ws.Range("L116").Value2 = 6037.8335
$ws.Range("M116").Value2 = -1219.5
$ws.Range("N116").Value2 = -12921.8335
$ws.Range("H125").Value2 = 7939931.5
$ws.Range("J125").Value2 = 9012657
$ws.Range("L125").Value2 = 81113913
$ws.Range("N125").Value2 = -81118833
$ws.Range("H132").Value2 = 20409762
$ws.Range("I132").Value2 = 21740702
$ws.Range("J132").Value2 = 1998.3334
$ws.Range("K132").Value2 = 65222106
$ws.Range("L132").Value2 = 5995.0002
$ws.Range("M132").Value2 = -65219576
$ws.Range("N132").Value2 = -11055.0002
$ws.Range("H137").Value2 = 3610.077
$ws.Range("I137").Value2 = 3585.611
$ws.Range("J137").Value2 = 3631.0476
$ws.Range("K137").Value2 = 10756.833
$ws.Range("L137").Value2 = 10893.1428
$ws.Range("M137").Value2 = -8206.832999999999
$ws.Range("N137").Value2 = -15993.1428
$ws.Range("H138").Value2 = 3302.0854
$ws.Range("I138").Value2 = 2186.7407
$ws.Range("K138").Value2 = 6560.222099999999
$ws.Range("M138").Value2 = -1420.222099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value2 = 7575614
$ws.Range("J45").Value2 = 10161.167
$ws.Range("L45").Value2 = 10161.167
$ws.Range("N45").Value2 = -10915.167
$ws.Range("H61").Value2 = 4080.8445
$ws.Range("I61").Value2 = 4142.927
$ws.Range("J61").Value2 = 3444.5
$ws.Range("K61").Value2 = 4142.927
$ws.Range("L61").Value2 = 3444.5
$ws.Range("M61").Value2 = -3930.927
$ws.Range("N61").Value2 = -3868.5
$ws.Range("H110").Value2 = 897001.5
$ws.Range("I110").Value2 = 958705.0600000001
$ws.Range("K110").Value2 = 958705.0600000001
$ws.Range("M110").Value2 = -956660.0600000001
$ws.Range("H136").Value2 = 4080.8445
$ws.Range("I136").Value2 = 4142.927
$ws.Range("J136").Value2 = 3444.5
$ws.Range("K136").Value2 = 12428.781
$ws.Range("L136").Value2 = 10333.5
$ws.Range("M136").Value2 = -9878.780999999999
$ws.Range("N136").Value2 = -15433.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 650
$ws.Range("J22").Value2 = 916.6667
$ws.Range("L22").Value2 = 916.6667
$ws.Range("N22").Value2 = -1262.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 2015.8
$ws.Range("I16").Value2 = 1438.9
$ws.Range("K16").Value2 = 1438.9
$ws.Range("M16").Value2 = -1151.9
$ws.Range("H82").Value2 = 32498.4
$ws.Range("I82").Value2 = 4999
$ws.Range("K82").Value2 = 4999
$ws.Range("M82").Value2 = -4638
$ws.Range("H85").Value2 = 32498.4
$ws.Range("I85").Value2 = 4999
$ws.Range("K85").Value2 = 4999
$ws.Range("M85").Value2 = -3751
$ws.Range("H106").Value2 = 29996.5
$ws.Range("J106").Value2 = 29996.5
$ws.Range("L106").Value2 = 29996.5
$ws.Range("N106").Value2 = -32520.5
$ws.Range("H113").Value2 = 2015.8
$ws.Range("I113").Value2 = 1438.9
$ws.Range("K113").Value2 = 1438.9
$ws.Range("M113").Value2 = 731.0999999999999
$ws.Range("H132").Value2 = 1840.7142
$ws.Range("J132").Value2 = 2577.1667
$ws.Range("L132").Value2 = 7731.500100000001
$ws.Range("N132").Value2 = -12791.5001
$ws.Range("H134").Value2 = 23662.717
$ws.Range("I134").Value2 = 27765.025
$ws.Range("J134").Value2 = 12234.857
$ws.Range("K134").Value2 = 83295.07500000001
$ws.Range("L134").Value2 = 36704.571
$ws.Range("M134").Value2 = -80760.07500000001
$ws.Range("N134").Value2 = -41774.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value2 = 15631304
$ws.Range("I56").Value2 = 15631304
$ws.Range("K56").Value2 = 15631304
$ws.Range("M56").Value2 = -15630774
$ws.Range("H75").Value2 = 707.6667
$ws.Range("J75").Value2 = 100
$ws.Range("L75").Value2 = 300
$ws.Range("N75").Value2 = -2296
$ws.Range("H78").Value2 = 707.6667
$ws.Range("J78").Value2 = 100
$ws.Range("L78").Value2 = 900
$ws.Range("N78").Value2 = -10884
$ws.Range("H92").Value2 = 609.7143
$ws.Range("I92").Value2 = 604.6
$ws.Range("J92").Value2 = 622.5
$ws.Range("K92").Value2 = 1813.8
$ws.Range("L92").Value2 = 1867.5
$ws.Range("M92").Value2 = -565.8000000000002
$ws.Range("N92").Value2 = -4363.5
$ws.Range("H122").Value2 = 886.2632
$ws.Range("I122").Value2 = 955.8570999999999
$ws.Range("J122").Value2 = 845.6667
$ws.Range("K122").Value2 = 8602.713899999999
$ws.Range("L122").Value2 = 7611.0003
$ws.Range("M122").Value2 = -6152.713899999999
$ws.Range("N122").Value2 = -12511.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value2 = 3416.5
$ws.Range("I55").Value2 = 2299.8
$ws.Range("J55").Value2 = 9000
$ws.Range("K55").Value2 = 2299.8
$ws.Range("L55").Value2 = 9000
$ws.Range("M55").Value2 = -1972.8
$ws.Range("N55").Value2 = -9654
$ws.Range("H58").Value2 = 17999
$ws.Range("J58").Value2 = 19999
$ws.Range("L58").Value2 = 19999
$ws.Range("N58").Value2 = -20553
$ws.Range("H122").Value2 = 63635.066
$ws.Range("I122").Value2 = 78789.75
$ws.Range("J122").Value2 = 3016.3333
$ws.Range("K122").Value2 = 236369.25
$ws.Range("L122").Value2 = 9048.999899999999
$ws.Range("M122").Value2 = -233919.25
$ws.Range("N122").Value2 = -13948.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 3000
$ws.Range("I22").Value2 = 0
$ws.Range("J22").Value2 = 3000
$ws.Range("K22").Value2 = 0
$ws.Range("L22").Value2 = 3000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value2 = -3590
$ws.Range("H27").Value2 = 3000
$ws.Range("I27").Value2 = 0
$ws.Range("J27").Value2 = 3000
$ws.Range("K27").Value2 = 0
$ws.Range("L27").Value2 = 3000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value2 = -3214
$ws.Range("H46").Value2 = 1455268.6
$ws.Range("J46").Value2 = 7677.077
$ws.Range("L46").Value2 = 7677.077
$ws.Range("N46").Value2 = -8053.077
$ws.Range("H93").Value2 = 13336912
$ws.Range("I93").Value2 = 19611396
$ws.Range("J93").Value2 = 3633.125
$ws.Range("K93").Value2 = 19611396
$ws.Range("L93").Value2 = 3633.125
$ws.Range("M93").Value2 = -19610148
$ws.Range("N93").Value2 = -6129.125
$ws.Range("H100").Value2 = 37406.734
$ws.Range("J100").Value2 = 146499.72
$ws.Range("L100").Value2 = 146499.72
$ws.Range("N100").Value2 = -147581.72

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value2 = 25011446
$ws.Range("I132").Value2 = 30313288
$ws.Range("J132").Value2 = 17044.715
$ws.Range("K132").Value2 = 90939864
$ws.Range("L132").Value2 = 51134.145
$ws.Range("M132").Value2 = -90937334
$ws.Range("N132").Value2 = -56194.145
